# Auto-generated edit script: update cryptos list values
# Preserves original cell style while forcing text storage so
# numeric-looking strings (e.g. '1.00', '94.170.86') are not
# silently reinterpreted by Excel's smart-parsing as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '94.170.86'
Set-TextValue $ws.Range('E2') '  +2.06%  '
Set-TextValue $ws.Range('D3') '3.070.15'
Set-TextValue $ws.Range('E3') '  -1.26%  '
Set-TextValue $ws.Range('E4') '  +0.04%  '
Set-TextValue $ws.Range('D5') '235.35'
Set-TextValue $ws.Range('E5') '  -1.00%  '
Set-TextValue $ws.Range('D6') '605.66'
Set-TextValue $ws.Range('E6') '  -1.67%  '
Set-TextValue $ws.Range('D7') '1.10'
Set-TextValue $ws.Range('E7') '  +0.07%  '
Set-TextValue $ws.Range('D8') '0.377'
Set-TextValue $ws.Range('E8') '  -3.55%  '
Set-TextValue $ws.Range('D9') '1.00'
Set-TextValue $ws.Range('E9') '  +0.10%  '
Set-TextValue $ws.Range('D10') '0.799'
Set-TextValue $ws.Range('E10') '  +8.09%  '
Set-TextValue $ws.Range('D11') '3.068.46'
Set-TextValue $ws.Range('E11') '  -1.22%  '
Set-TextValue $ws.Range('E12') '  -2.40%  '
Set-TextValue $ws.Range('D13') '93.796.45'
Set-TextValue $ws.Range('E13') '  +1.67%  '
Set-TextValue $ws.Range('D14') '0.0000239'
Set-TextValue $ws.Range('E14') '  -3.66%  '
Set-TextValue $ws.Range('D15') '33.57'
Set-TextValue $ws.Range('E15') '  -1.93%  '
Set-TextValue $ws.Range('D16') '5.29'
Set-TextValue $ws.Range('E16') '  -2.96%  '
Set-TextValue $ws.Range('D17') '3.644.63'
Set-TextValue $ws.Range('E17') '  -1.58%  '
Set-TextValue $ws.Range('D18') '3.067.10'
Set-TextValue $ws.Range('E18') '  -1.53%  '
Set-TextValue $ws.Range('D19') '3.53'
Set-TextValue $ws.Range('E19') '  -5.37%  '
Set-TextValue $ws.Range('D20') '14.24'
Set-TextValue $ws.Range('E20') '  -2.62%  '
Set-TextValue $ws.Range('D21') '5.71'
Set-TextValue $ws.Range('E21') '  -1.57%  '
Set-TextValue $ws.Range('D22') '439.72'
Set-TextValue $ws.Range('E22') '  -1.58%  '
Set-TextValue $ws.Range('D23') '8.81'
Set-TextValue $ws.Range('E23') '  -6.26%  '
Set-TextValue $ws.Range('D24') '0.0000187'
Set-TextValue $ws.Range('E24') '  -4.90%  '
Set-TextValue $ws.Range('D25') '8.38'
Set-TextValue $ws.Range('D26') '5.49'
Set-TextValue $ws.Range('E26') '  -5.08%  '
Set-TextValue $ws.Range('D27') '84.35'
Set-TextValue $ws.Range('E27') '  -2.72%  '
Set-TextValue $ws.Range('D28') '11.79'
Set-TextValue $ws.Range('E28') '  -0.20%  '
Set-TextValue $ws.Range('D29') '3.234.92'
Set-TextValue $ws.Range('E29') '  -1.29%  '
Set-TextValue $ws.Range('D30') '1.00'
Set-TextValue $ws.Range('D31') '0.247'
Set-TextValue $ws.Range('E31') '  +6.64%  '
Set-TextValue $ws.Range('D32') '0.177'
Set-TextValue $ws.Range('E32') '  +3.64%  '
Set-TextValue $ws.Range('D33') '0.122'
Set-TextValue $ws.Range('E33') '  -8.49%  '
Set-TextValue $ws.Range('D34') '1.00'
Set-TextValue $ws.Range('E34') '  +38.36%  '
Set-TextValue $ws.Range('D35') '8.82'
Set-TextValue $ws.Range('E35') '  -3.52%  '
Set-TextValue $ws.Range('D36') '7.37'
Set-TextValue $ws.Range('E36') '  -6.42%  '
Set-TextValue $ws.Range('E37') '  -4.61%  '
Set-TextValue $ws.Range('D38') '25.37'
Set-TextValue $ws.Range('E38') '  -3.08%  '
Set-TextValue $ws.Range('D39') '1.87'
Set-TextValue $ws.Range('E39') '  -1.75%  '
Set-TextValue $ws.Range('D40') '481.07'
Set-TextValue $ws.Range('E40') '  -0.39%  '
Set-TextValue $ws.Range('D41') '3.81'
Set-TextValue $ws.Range('E41') '  -1.21%  '
Set-TextValue $ws.Range('D42') '24.03'
Set-TextValue $ws.Range('E42') '  +0.74%  '
Set-TextValue $ws.Range('D43') '0.432'
Set-TextValue $ws.Range('E43') '  -0.30%  '
Set-TextValue $ws.Range('E44') '  -4.52%  '
Set-TextValue $ws.Range('E45') '  -0.02%  '
Set-TextValue $ws.Range('D46') '3.04'
Set-TextValue $ws.Range('E46') '  -7.66%  '
Set-TextValue $ws.Range('D47') '161.04'
Set-TextValue $ws.Range('E47') '  -0.55%  '
Set-TextValue $ws.Range('D48') '0.671'
Set-TextValue $ws.Range('E48') '  -3.04%  '
Set-TextValue $ws.Range('D49') '1.81'
Set-TextValue $ws.Range('E49') '  -4.35%  '
Set-TextValue $ws.Range('D50') '43.50'
Set-TextValue $ws.Range('E50') '  -1.05%  '
Set-TextValue $ws.Range('B51') 'FirstDigitalUSD'
Set-TextValue $ws.Range('C51') 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue $ws.Range('D51') '0.998'
Set-TextValue $ws.Range('E51') '  +0.12%  '
